$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B (old B/C shift to C/D)
$ws.Columns.Item(2).Insert()

# New column B content (header + stat query), matching width of column A
$ws.Columns.Item(2).ColumnWidth = 75

$ws.Cells.Item(1, 2).Value = "StatQuery"
$ws.Cells.Item(2, 2).Value = "MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE diag.disease_term IN['Lymphoma :: Stage 3']  OPTIONAL MATCH (f:file)-[*]->(c), (samp:sample)-[*]->(c) WITH DISTINCT c AS c, p, s, demo, diag, f, samp RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(samp)) as number_of_sample , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(s.clinical_study_designation)) as number_of_study"

# Apply the same wrap-text style used by A2 to the new B2 cell
$ws.Cells.Item(2, 2).WrapText = $true

# Update the selection to match the target state
$ws.Range("A2").Select()
